$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.058.20'
$ws.Range('D3').Value = '1.911.75'
$ws.Range('E3').Value = '  +0.49%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7905'
$ws.Range('E5').Value = '  +6.77%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '243.10'
$ws.Range('E6').Value = '  +0.10%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.000'
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('E8').Value = '  +3.53%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '26.40'
$ws.Range('E9').Value = '  +1.96%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06932'
$ws.Range('E10').Value = '  +0.41%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07988'
$ws.Range('E11').Value = '  -0.26%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.7490'
$ws.Range('E12').Value = '  -1.61%  '
$ws.Range('D13').Value = '1.913.52'
$ws.Range('E13').Value = '  +0.65%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.241'
$ws.Range('E14').Value = '  -0.01%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '93.50'
$ws.Range('E15').Value = '  +2.37%  '
$ws.Range('D16').Value = '30.071.57'
$ws.Range('E16').Value = '  +0.13%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '14.07'
$ws.Range('E17').Value = '  +0.20%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '5.945'
$ws.Range('E18').Value = '  -4.71%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '247.23'
$ws.Range('E19').Value = '  +4.03%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007798'
$ws.Range('E20').Value = '  +0.51%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.000'
$ws.Range('E21').Value = '  -0.02%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.000'
$ws.Range('E22').Value = '  -0.04%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.919'
$ws.Range('E23').Value = '  -1.89%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '170.02'
$ws.Range('E24').Value = '  +2.13%  '
$ws.Range('E25').Value = '  +0.13%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1392'
$ws.Range('E26').Value = '  +10.76%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.94'
$ws.Range('E27').Value = '  +0.61%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.048'
$ws.Range('E28').Value = '  +0.34%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.373'
$ws.Range('E29').Value = '  +1.66%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.525'
$ws.Range('E30').Value = '  -0.64%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.347'
$ws.Range('E31').Value = '  +1.02%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.05601'
$ws.Range('E32').Value = '  +5.17%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.118'
$ws.Range('E33').Value = '  +1.71%  '
$ws.Range('E34').Value = '  -2.30%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7379'
$ws.Range('E35').Value = '  -0.20%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.732'
$ws.Range('E36').Value = '  +0.13%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01935'
$ws.Range('E37').Value = '  -0.24%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.796'
$ws.Range('E38').Value = '  +0.14%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.186'
$ws.Range('E39').Value = '  -1.33%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.4450'
$ws.Range('E40').Value = '  +0.00%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '72.66'
$ws.Range('E41').Value = '  -0.59%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.000'
$ws.Range('E42').Value = '  +0.01%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.902'
$ws.Range('E43').Value = '  -3.16%  '
$ws.Range('E44').Value = '  -0.10%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '7.581'
$ws.Range('E45').Value = '  -0.61%  '
$ws.Range('E46').Value = '  -0.72%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.810'
$ws.Range('E47').Value = '  +0.21%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '989.03'
$ws.Range('E48').Value = '  +8.26%  '
$ws.Range('D49').Value = '2.066.04'
$ws.Range('E49').Value = '  +0.76%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '36.27'
$ws.Range('E50').Value = '  -1.20%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.501'
$ws.Range('E51').Value = '  +2.88%  '
